{"js": "// Replace each piece of text (the date line + the 25 multiplication problems)\n// with its updated value, matching the unified diff exactly. Every occurrence\n// of the \"old\" text in the document is unique, so a simple exact search +\n// Replace insertText is sufficient and keeps the original run formatting\n// (fonts, sizes, paragraph alignment, etc.) untouched.\nconst replacements = [\n  [\"2024-11-08 Friday\", \"2024-11-09 Saturday\"],\n  [\"283\u00d75=1415\", \"480\u00d73=1440\"],\n  [\"819\u00d72=1638\", \"168\u00d76=1008\"],\n  [\"907\u00d73=2721\", \"699\u00d74=2796\"],\n  [\"471\u00d76=2826\", \"209\u00d72=418\"],\n  [\"360\u00d75=1800\", \"996\u00d77=6972\"],\n  [\"365\u00d76=2190\", \"430\u00d75=2150\"],\n  [\"360\u00d72=720\", \"149\u00d78=1192\"],\n  [\"362\u00d73=1086\", \"895\u00d73=2685\"],\n  [\"397\u00d73=1191\", \"366\u00d79=3294\"],\n  [\"280\u00d74=1120\", \"517\u00d72=1034\"],\n  [\"164\u00d79=1476\", \"976\u00d79=8784\"],\n  [\"284\u00d75=1420\", \"671\u00d77=4697\"],\n  [\"782\u00d75=3910\", \"371\u00d77=2597\"],\n  [\"884\u00d77=6188\", \"559\u00d76=3354\"],\n  [\"994\u00d78=7952\", \"241\u00d73=723\"],\n  [\"835\u00d73=2505\", \"761\u00d74=3044\"],\n  [\"595\u00d74=2380\", \"321\u00d74=1284\"],\n  [\"733\u00d74=2932\", \"574\u00d79=5166\"],\n  [\"779\u00d77=5453\", \"723\u00d76=4338\"],\n  [\"191\u00d75=955\", \"292\u00d74=1168\"],\n  [\"423\u00d76=2538\", \"451\u00d78=3608\"],\n  [\"258\u00d74=1032\", \"464\u00d73=1392\"],\n  [\"513\u00d72=1026\", \"256\u00d76=1536\"],\n  [\"723\u00d79=6507\", \"644\u00d75=3220\"],\n  [\"746\u00d75=3730\", \"517\u00d72=1034\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // There is exactly one occurrence of each string in this document, but\n  // guard against duplicates anyway by replacing every hit.\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each piece of text (the date line + the 25 multiplication problems)\n# with its updated value, matching the unified diff exactly. Each \"old\" string\n# is unique in the document, so Find/Replace (wdReplaceAll) for each pair is\n# sufficient and preserves the original run formatting (fonts, sizes,\n# paragraph alignment, etc.).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-11-08 Friday\", \"2024-11-09 Saturday\"),\n    @(\"283\u00d75=1415\", \"480\u00d73=1440\"),\n    @(\"819\u00d72=1638\", \"168\u00d76=1008\"),\n    @(\"907\u00d73=2721\", \"699\u00d74=2796\"),\n    @(\"471\u00d76=2826\", \"209\u00d72=418\"),\n    @(\"360\u00d75=1800\", \"996\u00d77=6972\"),\n    @(\"365\u00d76=2190\", \"430\u00d75=2150\"),\n    @(\"360\u00d72=720\", \"149\u00d78=1192\"),\n    @(\"362\u00d73=1086\", \"895\u00d73=2685\"),\n    @(\"397\u00d73=1191\", \"366\u00d79=3294\"),\n    @(\"280\u00d74=1120\", \"517\u00d72=1034\"),\n    @(\"164\u00d79=1476\", \"976\u00d79=8784\"),\n    @(\"284\u00d75=1420\", \"671\u00d77=4697\"),\n    @(\"782\u00d75=3910\", \"371\u00d77=2597\"),\n    @(\"884\u00d77=6188\", \"559\u00d76=3354\"),\n    @(\"994\u00d78=7952\", \"241\u00d73=723\"),\n    @(\"835\u00d73=2505\", \"761\u00d74=3044\"),\n    @(\"595\u00d74=2380\", \"321\u00d74=1284\"),\n    @(\"733\u00d74=2932\", \"574\u00d79=5166\"),\n    @(\"779\u00d77=5453\", \"723\u00d76=4338\"),\n    @(\"191\u00d75=955\", \"292\u00d74=1168\"),\n    @(\"423\u00d76=2538\", \"451\u00d78=3608\"),\n    @(\"258\u00d74=1032\", \"464\u00d73=1392\"),\n    @(\"513\u00d72=1026\", \"256\u00d76=1536\"),\n    @(\"723\u00d79=6507\", \"644\u00d75=3220\"),\n    @(\"746\u00d75=3730\", \"517\u00d72=1034\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
